# Fix navigation drawer selectable
# - Move the "Navigation Drawer / Playlist highlight seleted item / FIX" row
#   (row 37) up to directly follow the other "Navigation Drawer" FIX rows
#   (now row 29), shifting the UPDATE rows that used to be 29-36 down to 30-37.
# - Re-normalize row 25's leftover blank D/F/G cells to match its neighbours.
# - Update the saved window/selection/sort-range metadata to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Stage old row 37 (values + full formatting) in a scratch row so it
#    survives while rows 29-36 shift down into 30-37.
# ---------------------------------------------------------------------
$ws.Range("A37:G37").Copy()
$ws.Range("A200:G200").PasteSpecial(-4122)
$ws.Range("A200:G200").Value = $ws.Range("A37:G37").Value()

# ---------------------------------------------------------------------
# 2) Shift rows 29-36 down into 30-37 (values + formatting), bottom-up so
#    we never overwrite a source row before it has been read.
# ---------------------------------------------------------------------
for ($r = 36; $r -ge 29; $r--) {
    $src = $r
    $dst = $r + 1
    $ws.Range("A$($src):G$($src)").Copy()
    $ws.Range("A$($dst):G$($dst)").PasteSpecial(-4122)
    $ws.Range("A$($dst):G$($dst)").Value = $ws.Range("A$($src):G$($src)").Value()
}

# ---------------------------------------------------------------------
# 3) Drop the staged old-row-37 content into row 29 (keeps the blank
#    D29:G29 formatted like D28:G28 / the rest of the A:G block).
# ---------------------------------------------------------------------
$ws.Range("A200:G200").Copy()
$ws.Range("A29:G29").PasteSpecial(-4122)
$ws.Range("A29:G29").Value = $ws.Range("A200:G200").Value()

# Clear the scratch row entirely.
$ws.Range("A200:G200").Clear()

# D29:G29 must stay truly blank (no leftover values from the old row 37,
# which only had A:C populated).
$ws.Range("D29:G29").ClearContents()

# Rows 33 and 34 (content shifted from the old un-annotated rows 32/33)
# never had D/F/G cells at all - drop the blank placeholders the format
# copy introduced so they fully disappear again (matches their donor rows,
# which only ever had A:C + E populated). E33/E34 keep their "HaiNNT" value.
$ws.Range("D33").Clear()
$ws.Range("F33").Clear()
$ws.Range("G33").Clear()
$ws.Range("D34").Clear()
$ws.Range("F34").Clear()
$ws.Range("G34").Clear()

# ---------------------------------------------------------------------
# 4) Row 25: normalize format of A/B/C/E to match row 24's style and drop
#    the now-unused blank D/F/G cells entirely.
# ---------------------------------------------------------------------
$ws.Range("A24:G24").Copy()
$ws.Range("A25:G25").PasteSpecial(-4122)
$ws.Range("D25").Clear()
$ws.Range("F25").Clear()
$ws.Range("G25").Clear()

